# Added variables needed for changelogs.
# Adds a new "changes" variable (Changes Made / LongText) to the Variables
# sheet, tags it (and the existing "rev" row) with the Assets category in
# addition to Core/Customer/Catering, and mirrors the new column into each
# of the four per-category SQL-table sheets.

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCore      = $wb.Worksheets.Item("Core SQL Table")
$wsCustomer  = $wb.Worksheets.Item("Customer SQL Table")
$wsCatering  = $wb.Worksheets.Item("Catering SQL Table")
$wsAssets    = $wb.Worksheets.Item("Assets SQL Table")

# --- Variables sheet ---------------------------------------------------
# New row 29: the "changes" variable itself (write first so the shared
# string table picks up "changes"/"Changes Made" before the now-longer
# "Core/Customer/Catering/Assets" string used below).
$wsVariables.Cells.Item(29, 1).Value = "changes"
$wsVariables.Cells.Item(29, 2).Value = "Changes Made"
$wsVariables.Cells.Item(29, 3).Value = "LongText"
$wsVariables.Cells.Item(29, 4).Value = "Core/Customer/Catering/Assets"

# rev (row 26) now also belongs to the Assets database/category.
$wsVariables.Cells.Item(26, 4).Value = "Core/Customer/Catering/Assets"

# Widen column D so the longer category string keeps fitting.
$wsVariables.Columns.Item(4).ColumnWidth = 29.25

# --- Per-category SQL table sheets: append the new "changes" column ----
$wsCore.Cells.Item(1, 17).Value = "changes"      # Column Q
$wsCustomer.Cells.Item(1, 8).Value = "changes"   # Column H
$wsCatering.Cells.Item(1, 4).Value = "changes"   # Column D
$wsAssets.Cells.Item(1, 4).Value = "changes"     # Column D

# --- Selections (set non-active sheets first, Variables last so it stays
#     the active tab, matching the saved workbook view state) ----------
$wsCore.Range("Q1").Select()
$wsCustomer.Range("H1").Select()
$wsCatering.Range("D1").Select()
$wsAssets.Range("D1").Select()
$wsVariables.Range("D30").Select()
